$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.057.81'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.26%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.173.88'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.40%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.09%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.57%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.40'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +7.29%  '

# Row 7
$ws.Range("E7").Value = '  +0.05%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.169.77'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.29%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.532'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.92%  '

# Row 10
$ws.Range("E10").Value = '  +6.80%  '

# Row 11
$ws.Range("E11").Value = '  +0.22%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.504'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.27%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000271'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +19.01%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.62'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.26%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.689.73'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.34%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.104.90'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.27%  '

# Row 17
$ws.Range("B17").Value = 'Polkadot'
$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.19'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +6.80%  '

# Row 18
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.171.89'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.27%  '

# Row 19
$ws.Range("E19").Value = '  +1.57%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '512.22'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.71%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.90'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.10%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.728'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +7.20%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '15.46'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.76%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.84'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.51%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.18'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.47%  '

# Row 26
$ws.Range("E26").Value = '  +0.10%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.15'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +13.99%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.94'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.54%  '

# Row 29
$ws.Range("E29").Value = '  +9.26%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '27.88'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +7.09%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.81'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +15.70%  '

# Row 32
$ws.Range("E32").Value = '  -0.04%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.20'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.36%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.33'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +12.66%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.61'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +7.02%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '55.76'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.39%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0903'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +11.62%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '476.04'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.29%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.11'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +13.66%  '

# Row 40
$ws.Range("E40").Value = '  +3.73%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.70'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.27%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.078.80'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.86%  '

# Row 43
$ws.Range("E43").Value = '  +2.15%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.46'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +10.88%  '

# Row 45
$ws.Range("E45").Value = '  +6.78%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '29.37'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.99%  '

# Row 47
$ws.Range("E47").Value = '  +19.98%  '

# Row 49
$ws.Range("E49").Value = '  +1.47%  '

# Row 50
$ws.Range("E50").Value = '  +9.21%  '

# Row 51
$ws.Range("E51").Value = '  +2.04%  '
